# "rolling back hot fix due to extra errors"
#
# This restores the workbook to its pre-hotfix state:
#   - the sheet "Project - Contact" is renamed back to "Contact"
#   - the previously-selected tab ("Project - Contact"/"Contact") is no
#     longer the active tab
#   - the "Donor organism" sheet becomes the active tab again, with its
#     selection restored to cell AD28

$wb = $excel.ActiveWorkbook

# Rename "Project - Contact" back to "Contact"
$wsContact = $wb.Worksheets.Item("Project - Contact")
$wsContact.Name = "Contact"

# Re-activate "Donor organism" as the selected sheet, restoring its
# previous cell selection (this also moves tabSelected/activeTab off of
# the Contact sheet and onto Donor organism)
$wsDonor = $wb.Worksheets.Item("Donor organism")
$wsDonor.Activate()
$wsDonor.Range("AD28").Select()
